# Informação de instrutor por projeto
# Rebuilds the "Instrutor x Projeto" consolidated table: the DD2 / IdearTec
# columns are split into per-wave (Onda1/2/3) columns, a new IT1 column is
# added, PROG_2 is renumbered/resequenced (a new PROG_12 appears), and the
# ROB_8 / ROB_9 rows are dropped - so the table shrinks from A1:F23 to A1:J20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- header row (A1:J1) -----------------------------------------------
$headers = New-Object 'object[,]' 1,10
$headers[0,0] = "Instrutor"
$headers[0,1] = "DD1"
$headers[0,2] = "DD2_Onda1"
$headers[0,3] = "DD2_Onda2"
$headers[0,4] = "DD2_Onda3"
$headers[0,5] = "IT1"
$headers[0,6] = "IT2_Onda1"
$headers[0,7] = "IT2_Onda2"
$headers[0,8] = "IT2_Onda3"
$headers[0,9] = "Total"
$ws.Range("A1:J1").Value = $headers

# Extend the existing bold/border/center header style (previously only on
# A1:F1) across the four newly-added header cells (G1:J1) by copying the
# format from the already-styled F1 cell.
$ws.Range("F1").Copy()
$ws.Range("G1:J1").PasteSpecial(-4122)

# ---- data rows (A2:J20) -------------------------------------------------
$data = New-Object 'object[,]' 19,10
$rows = @(
    @("PROG_1",0,1,1,1,2,4,4,1,14),
    @("PROG_10",1,4,3,1,1,1,1,2,14),
    @("PROG_11",1,1,1,1,1,2,2,2,11),
    @("PROG_12",2,2,3,3,0,0,1,3,14),
    @("PROG_2",0,0,1,1,2,5,2,3,14),
    @("PROG_3",1,1,3,2,1,2,1,3,14),
    @("PROG_4",0,1,1,2,3,1,0,6,14),
    @("PROG_5",0,2,1,2,2,2,4,1,14),
    @("PROG_6",0,1,2,3,2,2,4,0,14),
    @("PROG_7",1,2,2,2,1,2,2,2,14),
    @("PROG_8",1,3,3,2,3,1,0,1,14),
    @("PROG_9",1,3,0,1,2,2,3,2,14),
    @("ROB_1",0,1,1,1,0,4,3,2,12),
    @("ROB_2",0,1,3,1,0,1,2,3,11),
    @("ROB_3",0,0,4,1,0,0,2,3,10),
    @("ROB_4",0,5,0,3,0,1,0,1,10),
    @("ROB_5",0,1,4,1,0,1,2,1,10),
    @("ROB_6",0,4,0,3,0,1,1,1,10),
    @("ROB_7",0,2,2,4,0,2,0,0,10)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $data[$r, $c] = $row[$c]
    }
}

$ws.Range("A2:J20").Value = $data

# ---- drop the old trailing rows (previously ROB_7..ROB_9, now unused) ---
$ws.Range("A21:F23").ClearContents()
